$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 77; this shifts existing rows 77..157 down to 78..158
$ws.Rows.Item(77).Insert()

# Fill in the new row 77 with the new record's data
$ws.Cells.Item(77, 1).Value = 11
$ws.Cells.Item(77, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(77, 3).Value = "Bíobío"
$ws.Cells.Item(77, 4).Value = 45072
$ws.Cells.Item(77, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(77, 5).Value = 8
$ws.Cells.Item(77, 6).Value = 100112001
$ws.Cells.Item(77, 7).Value = "Berenjena"
$ws.Cells.Item(77, 8).Value = "Sin especificar"
$ws.Cells.Item(77, 9).Value = "Primera"
$ws.Cells.Item(77, 10).Value = 350
$ws.Cells.Item(77, 11).Value = 7000
$ws.Cells.Item(77, 12).Value = 7500
$ws.Cells.Item(77, 13).Value = 7286
$ws.Cells.Item(77, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(77, 15).Value = "Región Metropolitana"
$ws.Cells.Item(77, 16).Value = 121
$ws.Cells.Item(77, 17).Value = 60
$ws.Cells.Item(77, 18).Value = "Hortaliza"
